# Auto-generated PowerShell COM-interop script
# Applies numeric cell updates to the profit-calculation sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H..N) per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 934.6667
$ws.Range("I29").Value = 30
$ws.Range("J29").Value = 1387
$ws.Range("K29").Value = 90
$ws.Range("L29").Value = 4161
$ws.Range("M29").Value = 191
$ws.Range("N29").Value = -4723
$ws.Range("H32").Value = 1168.4286
$ws.Range("J32").Value = 1115.8
$ws.Range("L32").Value = 1115.8
$ws.Range("N32").Value = -1767.8
$ws.Range("H62").Value = 100045770
$ws.Range("I62").Value = 250001810
$ws.Range("K62").Value = 250001810
$ws.Range("M62").Value = -250001186
$ws.Range("H65").Value = 100045770
$ws.Range("I65").Value = 250001810
$ws.Range("K65").Value = 1250009050
$ws.Range("M65").Value = -1250005930
$ws.Range("H76").Value = 117528.5
$ws.Range("I76").Value = 128448.336
$ws.Range("J76").Value = 19250
$ws.Range("K76").Value = 128448.336
$ws.Range("L76").Value = 19250
$ws.Range("M76").Value = -128133.336
$ws.Range("N76").Value = -19880
$ws.Range("H79").Value = 117528.5
$ws.Range("I79").Value = 128448.336
$ws.Range("J79").Value = 19250
$ws.Range("K79").Value = 128448.336
$ws.Range("L79").Value = 19250
$ws.Range("M79").Value = -127356.336
$ws.Range("N79").Value = -21434
$ws.Range("H98").Value = 5540.9766
$ws.Range("I98").Value = 5479.3228
$ws.Range("K98").Value = 5479.3228
$ws.Range("M98").Value = -3981.3228
$ws.Range("H106").Value = 3599.5
$ws.Range("I106").Value = 3599.5
$ws.Range("K106").Value = 3599.5
$ws.Range("M106").Value = -2968.5
$ws.Range("H107").Value = 18751818
$ws.Range("I107").Value = 7355021
$ws.Range("K107").Value = 7355021
$ws.Range("M107").Value = -7353101
$ws.Range("H111").Value = 7816892
$ws.Range("J111").Value = 2481.3333
$ws.Range("L111").Value = 7443.999899999999
$ws.Range("N111").Value = -13577.9999
$ws.Range("H116").Value = 15629774
$ws.Range("I116").Value = 41668830
$ws.Range("J116").Value = 6339.2
$ws.Range("K116").Value = 41668830
$ws.Range("L116").Value = 6339.2
$ws.Range("M116").Value = -41665388
$ws.Range("N116").Value = -13223.2
$ws.Range("H122").Value = 5540.9766
$ws.Range("I122").Value = 5479.3228
$ws.Range("K122").Value = 16437.9684
$ws.Range("M122").Value = -13987.9684
$ws.Range("H132").Value = 1677.8125
$ws.Range("I132").Value = 1663.9131
$ws.Range("K132").Value = 4991.7393
$ws.Range("M132").Value = -2461.7393
$ws.Range("H137").Value = 2702
$ws.Range("J137").Value = 2399.9
$ws.Range("L137").Value = 7199.700000000001
$ws.Range("N137").Value = -12299.7
$ws.Range("H138").Value = 5151.3
$ws.Range("I138").Value = 1673.4445
$ws.Range("J138").Value = 6641.8096
$ws.Range("K138").Value = 5020.333500000001
$ws.Range("L138").Value = 19925.4288
$ws.Range("M138").Value = 119.6664999999994
$ws.Range("N138").Value = -30205.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 848
$ws.Range("J16").Value = 475
$ws.Range("L16").Value = 475
$ws.Range("N16").Value = -1049
$ws.Range("H45").Value = 3841.158
$ws.Range("I45").Value = 1226.75
$ws.Range("J45").Value = 5742.5454
$ws.Range("K45").Value = 1226.75
$ws.Range("L45").Value = 5742.5454
$ws.Range("M45").Value = -849.75
$ws.Range("N45").Value = -6496.5454
$ws.Range("H61").Value = 7514.263
$ws.Range("I61").Value = 4489.591
$ws.Range("K61").Value = 4489.591
$ws.Range("M61").Value = -4277.591
$ws.Range("H74").Value = 48712.227
$ws.Range("I74").Value = 68911.60000000001
$ws.Range("J74").Value = 5427.857
$ws.Range("K74").Value = 68911.60000000001
$ws.Range("L74").Value = 5427.857
$ws.Range("M74").Value = -68037.60000000001
$ws.Range("N74").Value = -7175.857
$ws.Range("H77").Value = 48712.227
$ws.Range("I77").Value = 68911.60000000001
$ws.Range("J77").Value = 5427.857
$ws.Range("K77").Value = 344558
$ws.Range("L77").Value = 27139.285
$ws.Range("M77").Value = -340190
$ws.Range("N77").Value = -35875.285
$ws.Range("H110").Value = 16668435
$ws.Range("I110").Value = 1797.7222
$ws.Range("K110").Value = 1797.7222
$ws.Range("M110").Value = 247.2778000000001
$ws.Range("H132").Value = 7472.1577
$ws.Range("I132").Value = 6486.619
$ws.Range("K132").Value = 19459.857
$ws.Range("M132").Value = -16929.857
$ws.Range("H136").Value = 7514.263
$ws.Range("I136").Value = 4489.591
$ws.Range("K136").Value = 13468.773
$ws.Range("M136").Value = -10918.773

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 38461864
$ws.Range("J80").Value = 273.8889
$ws.Range("L80").Value = 273.8889
$ws.Range("N80").Value = -2269.8889
$ws.Range("H83").Value = 38461864
$ws.Range("J83").Value = 273.8889
$ws.Range("L83").Value = 1369.4445
$ws.Range("N83").Value = -11353.4445
$ws.Range("H99").Value = 30306012
$ws.Range("I99").Value = 4450
$ws.Range("K99").Value = 4450
$ws.Range("M99").Value = -2952

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10017.061
$ws.Range("I31").Value = 4996.727
$ws.Range("J31").Value = 12527.228
$ws.Range("K31").Value = 4996.727
$ws.Range("L31").Value = 12527.228
$ws.Range("M31").Value = -4701.727
$ws.Range("N31").Value = -13117.228
$ws.Range("H34").Value = 10017.061
$ws.Range("I34").Value = 4996.727
$ws.Range("J34").Value = 12527.228
$ws.Range("K34").Value = 4996.727
$ws.Range("L34").Value = 12527.228
$ws.Range("M34").Value = -4794.727
$ws.Range("N34").Value = -12931.228
$ws.Range("H58").Value = 11117009
$ws.Range("I58").Value = 29413480
$ws.Range("J58").Value = 8436.964
$ws.Range("K58").Value = 29413480
$ws.Range("L58").Value = 8436.964
$ws.Range("M58").Value = -29413277
$ws.Range("N58").Value = -8842.964
$ws.Range("H99").Value = 5528.4546
$ws.Range("I99").Value = 2559.8
$ws.Range("J99").Value = 8002.3335
$ws.Range("K99").Value = 2559.8
$ws.Range("L99").Value = 8002.3335
$ws.Range("M99").Value = -1061.8
$ws.Range("N99").Value = -10998.3335
$ws.Range("H126").Value = 5528.4546
$ws.Range("I126").Value = 2559.8
$ws.Range("J126").Value = 8002.3335
$ws.Range("K126").Value = 7679.400000000001
$ws.Range("L126").Value = 24007.0005
$ws.Range("M126").Value = -5209.400000000001
$ws.Range("N126").Value = -28947.0005
$ws.Range("H136").Value = 11117009
$ws.Range("I136").Value = 29413480
$ws.Range("J136").Value = 8436.964
$ws.Range("K136").Value = 88240440
$ws.Range("L136").Value = 25310.892
$ws.Range("M136").Value = -88237890
$ws.Range("N136").Value = -30410.892

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 223431.89
$ws.Range("I2").Value = 148
$ws.Range("J2").Value = 502536.75
$ws.Range("K2").Value = 888
$ws.Range("L2").Value = 3015220.5
$ws.Range("M2").Value = -775
$ws.Range("N2").Value = -3015446.5
$ws.Range("H87").Value = 71438430
$ws.Range("I87").Value = 333335330
$ws.Range("K87").Value = 1000005990
$ws.Range("M87").Value = -1000004742
$ws.Range("H90").Value = 71438430
$ws.Range("I90").Value = 333335330
$ws.Range("K90").Value = 3000017970
$ws.Range("M90").Value = -3000011730
$ws.Range("H128").Value = 310000
$ws.Range("I128").Value = 310000
$ws.Range("K128").Value = 930000
$ws.Range("M128").Value = -925020

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5955.222
$ws.Range("I126").Value = 1999
$ws.Range("J126").Value = 6449.75
$ws.Range("K126").Value = 5997
$ws.Range("L126").Value = 19349.25
$ws.Range("M126").Value = -3527
$ws.Range("N126").Value = -24289.25
$ws.Range("H132").Value = 5551.625
$ws.Range("I132").Value = 1597.3334
$ws.Range("K132").Value = 4792.0002
$ws.Range("M132").Value = -2262.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 9615.679
$ws.Range("I136").Value = 4018.4614
$ws.Range("J136").Value = 14466.6
$ws.Range("K136").Value = 12055.3842
$ws.Range("L136").Value = 43399.8
$ws.Range("M136").Value = -9505.3842
$ws.Range("N136").Value = -48499.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16189.167
$ws.Range("J41").Value = 16189.167
$ws.Range("L41").Value = 16189.167
$ws.Range("N41").Value = -16969.167
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H62").Value = 6125
$ws.Range("J62").Value = 9899
$ws.Range("L62").Value = 9899
$ws.Range("N62").Value = -11147
$ws.Range("H65").Value = 6125
$ws.Range("J65").Value = 9899
$ws.Range("L65").Value = 49495
$ws.Range("N65").Value = -55735
$ws.Range("H81").Value = 16672877
$ws.Range("J81").Value = 100030000
$ws.Range("L81").Value = 200060000
$ws.Range("N81").Value = -200062122
$ws.Range("H84").Value = 16672877
$ws.Range("J84").Value = 100030000
$ws.Range("L84").Value = 1000300000
$ws.Range("N84").Value = -1000310608
$ws.Range("H107").Value = 1149.3636
$ws.Range("I107").Value = 1124.3334
$ws.Range("J107").Value = 1179.4
$ws.Range("K107").Value = 3373.0002
$ws.Range("L107").Value = 3538.2
$ws.Range("M107").Value = -1453.0002
$ws.Range("N107").Value = -7378.200000000001
$ws.Range("H109").Value = 59190
$ws.Range("J109").Value = 59190
$ws.Range("L109").Value = 59190
$ws.Range("N109").Value = -61964
$ws.Range("H122").Value = 26530414
$ws.Range("I122").Value = 42002850
$ws.Range("J122").Value = 6242.7144
$ws.Range("K122").Value = 126008550
$ws.Range("L122").Value = 18728.1432
$ws.Range("M122").Value = -126006100
$ws.Range("N122").Value = -23628.1432
$ws.Range("H126").Value = 2900.6072
$ws.Range("I126").Value = 2233.6667
$ws.Range("K126").Value = 6701.000100000001
$ws.Range("M126").Value = -4231.000100000001
